# Resize the 10 data columns on the "Data" sheet.
# Target stored OOXML <col> widths (per the diff):
#   col 1:  18 -> 15
#   col 2:  18 -> 15
#   col 3:  8.4 -> 7
#   col 4:  9.6 -> 8
#   col 5:  9.6 -> 8
#   col 6:  8.4 -> 7
#   col 7:  6 -> 5
#   col 8:  15.6 -> 13
#   col 9:  7.2 -> 6
#   col 10: 60 -> 50
#
# Excel's COM ColumnWidth property and the width stored in the XML differ by
# a constant offset (5/6, i.e. the built-in cell padding) for this
# Calibri 11 workbook, so we subtract that offset before assigning.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offset = 5 / 6

$targetWidths = @(15, 15, 7, 8, 8, 7, 5, 13, 6, 50)

for ($i = 0; $i -lt $targetWidths.Length; $i++) {
    $col = $i + 1
    $ws.Columns.Item($col).ColumnWidth = $targetWidths[$i] - $offset
}
